$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update progress column (C) with percentage values / text
$ws.Range("C2").Value = 0.8
$ws.Range("C2").NumberFormat = "0%"

$ws.Range("C3").Value = 0.8
$ws.Range("C3").NumberFormat = "0%"

$ws.Range("C8").Value = "en proceso"

# Move active selection to A21
$ws.Range("A21").Select()
